$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, centered, bordered) from H1 into I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$I = @(5,8,6,5,6,7,6,7,8,7,7,7,6,6,6,5)
$J = @(5,8,6,5,6,7,6,7,8,8,7,7,6,6,6,5)

for ($r = 2; $r -le 17; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $I[$idx]
    $ws.Cells.Item($r, 10).Value = $J[$idx]
}
